# Result tables updated after AutLvls changed
#
# The underlying data-generation run was repeated with updated AutLvls
# classification, which shifted a handful of the "overall"/"HC" aggregate
# percentages on the AutLvl==3 rows (row 6, TC06_* test case) in four of
# the result sheets: EmailsAllow_Observed_rel, HandsOffAllow_Observed_rel,
# Observed_Instr_rel and Observed_Rep_rel. Column C ("HC" overall) and
# column O ("overall") are the two columns that moved on each of those
# sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "EmailsAllow_Observed_rel"   = @{ "C6" = 95.63; "O6" = 95.24 }
    "HandsOffAllow_Observed_rel" = @{ "C6" = 91.67; "O6" = 90.48 }
    "Observed_Instr_rel"         = @{ "C6" = 90.87; "O6" = 95.24 }
    "Observed_Rep_rel"           = @{ "C6" = 95.24; "O6" = 100  }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellValues = $updates[$sheetName]
    foreach ($addr in $cellValues.Keys) {
        $ws.Range($addr).Value = $cellValues[$addr]
    }
}
